# Chart of Accounts update:
#  - Insert a new "Retained Earnings" (Equity) row right after "Tony's Equity"
#  - Append a new "Net Income" (Profit) row at the end
#  - Add two new columns: Column1 (Permanent/Temporary) and Natural_Balance (DEBIT/CREDIT)
#  - Clear the leftover manual highlighting that used to sit on the old rows 65/66
#  - Remove the duplicate-values conditional formatting rule on column E
#  - Refresh the table (ListObject) range/columns and the active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("chart_of_accounts")
$ws.Activate()

# --- 1. Insert the new "Retained Earnings" row after row 20 (GL 300101) ---
$ws.Rows.Item(21).Insert()
$ws.Cells.Item(21, 1).Value = 300102
$ws.Cells.Item(21, 2).Value = "Retained Earnings"
$ws.Cells.Item(21, 3).Value = "Equity"
$ws.Cells.Item(21, 4).Value = "Equity"

# --- 2. Append the new "Net Income" row at the bottom (now row 74) ---
$ws.Cells.Item(74, 1).Value = 700001
$ws.Cells.Item(74, 2).Value = "Net Income"
$ws.Cells.Item(74, 3).Value = "Profit"
$ws.Cells.Item(74, 4).Value = "Profit"

# --- 3. Grow the table to A1:D74, then add the two new columns ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:D74"))

$lo.ListColumns.Add() | Out-Null
$lo.ListColumns.Add() | Out-Null
$ws.Cells.Item(1, 5).Value = "Column1"
$ws.Cells.Item(1, 6).Value = "Natural_Balance"

# --- 4. Fill in Column1 / Natural_Balance for every data row based on Account_Type ---
$map = @{
    "Asset"      = @("Permanent", "DEBIT")
    "Liability"  = @("Permanent", "CREDIT")
    "Equity"     = @("Permanent", "CREDIT")
    "Revenue"    = @("Temporary", "CREDIT")
    "Deduction"  = @("Temporary", "DEBIT")
    "Expense"    = @("Temporary", "DEBIT")
    "Profit"     = @("Temporary", "CREDIT")
}

for ($r = 2; $r -le 74; $r++) {
    $atype = $ws.Cells.Item($r, 4).Value()
    $pair = $map[$atype]
    $ws.Cells.Item($r, 5).Value = $pair[0]
    $ws.Cells.Item($r, 6).Value = $pair[1]
}

# --- 5. Column F width, roughly matching the rest of the bestFit columns ---
$ws.Columns.Item(6).ColumnWidth = 17

# --- 6. Clear the old manual highlight styling that lived on rows 65/66          ---
#     (those rows are now 66/67 after the row-21 insert shifted everything down)
$ws.Range("B66:C67").ClearFormats()

# --- 7. Remove the "duplicate values" conditional formatting rule on column E ---
$cf = $ws.Range("E53:E72").FormatConditions
while ($cf.Count() -gt 0) {
    $cf.Item(1).Delete()
}

# --- 8. Restore a plausible active selection/view ---
$ws.Range("J41").Select()
